{"js": "// The author's edit centers the value cell (\"<0> Weeks\") that sits next to\n// the \"Age of birds at assessment\" label in the Management table \u2014 every\n// other \"<0>\" value cell in that table is already center-aligned; this one\n// was the odd one out, so bring it in line with its siblings.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet targetCell = null;\n\nfor (const table of tables.items) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (const row of rows.items) {\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    // Load the text of every cell in the row so we can find the label.\n    for (const cell of cells.items) {\n      cell.body.load(\"text\");\n    }\n    await context.sync();\n\n    const labelIndex = cells.items.findIndex(\n      (cell) => cell.body.text.trim() === \"Age of birds at assessment\"\n    );\n\n    if (labelIndex !== -1 && labelIndex + 1 < cells.items.length) {\n      targetCell = cells.items[labelIndex + 1];\n      break;\n    }\n  }\n\n  if (targetCell) break;\n}\n\nif (!targetCell) {\n  throw new Error(\"Could not locate the 'Age of birds at assessment' value cell.\");\n}\n\nconst paragraphs = targetCell.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Center every paragraph in that cell (there is exactly one here).\nfor (const paragraph of paragraphs.items) {\n  paragraph.alignment = Word.Alignment.centered;\n}\nawait context.sync();\n", "ps1": "# The author's edit centers the value cell (\"<0> Weeks\") that sits next to\n# the \"Age of birds at assessment\" label in the Management table - every\n# other \"<0>\" value cell in that table is already center-aligned; this one\n# was the odd one out, so bring it in line with its siblings.\n\n$wdAlignParagraphCenter = 1\n\n$d = $word.ActiveDocument\n\n$targetCell = $null\n\nforeach ($table in $d.Tables) {\n    foreach ($row in $table.Rows) {\n        $labelCell = $null\n        foreach ($cell in $row.Cells) {\n            $cellText = $cell.Range.Text.Trim().TrimEnd([char]7).TrimEnd([char]13)\n            if ($cellText -eq \"Age of birds at assessment\") {\n                $labelCell = $cell\n            }\n        }\n        if ($labelCell -ne $null) {\n            $targetCell = $row.Cells($labelCell.ColumnIndex + 1)\n            break\n        }\n    }\n    if ($targetCell -ne $null) {\n        break\n    }\n}\n\nif ($targetCell -eq $null) {\n    throw \"Could not locate the 'Age of birds at assessment' value cell.\"\n}\n\nforeach ($para in $targetCell.Range.Paragraphs) {\n    $para.Alignment = $wdAlignParagraphCenter\n}\n"}
